$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.049.06"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.58%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.616.37"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.71%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.02%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.526"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.615.52"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.165"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +14.32%  "
$ws.Range("E11").Value = "  +0.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.347"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.85%  "
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.097.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.59"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.98%  "
$ws.Range("E16").Value = "  +6.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "70.946.92"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.62%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.621.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "379.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.15"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.91%  "
$ws.Range("B24").Value = "NEARProtocol"
$ws.Range("C24").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.55%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.86"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.733.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0952"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "529.01"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.02"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.04%  "
$ws.Range("E33").Value = "  +3.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.83"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.47%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "165.42"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.118"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.14"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.97"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.55%  "
$ws.Range("E41").Value = "  +3.04%  "
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.89%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.01"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.331"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "154.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.36%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.530"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.01%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0264"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.52%  "
$ws.Range("E51").Value = "  +5.27%  "
